$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.960.10"
$ws.Range("E2").Value = "  -1.60%  "

$ws.Range("D3").Value = "3.397.83"
$ws.Range("E3").Value = "  -1.84%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "572.87"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.37%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "142.73"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.84%  "

$ws.Range("E7").Value = "  +0.09%  "

$ws.Range("D8").Value = "3.397.30"
$ws.Range("E8").Value = "  -1.86%  "

$ws.Range("E9").Value = "  +0.01%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.54"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.01%  "

$ws.Range("E11").Value = "  -0.39%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.394"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.41%  "

$ws.Range("D13").Value = "3.980.03"
$ws.Range("E13").Value = "  -1.69%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "28.27"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.25%  "

$ws.Range("E15").Value = "  +0.56%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000172"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.12%  "

$ws.Range("D17").Value = "3.401.91"
$ws.Range("E17").Value = "  -1.70%  "

$ws.Range("D18").Value = "61.099.82"
$ws.Range("E18").Value = "  -1.45%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.28"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.88%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.19"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.57%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "9.13"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.54%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "389.37"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.01%  "

$ws.Range("E23").Value = "  -0.65%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "73.19"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.96%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.996"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.71%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0000120"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.50%  "

$ws.Range("D27").Value = "3.534.70"
$ws.Range("E27").Value = "  -1.56%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.179"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.19%  "

$ws.Range("E29").Value = "  -0.20%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.41"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -5.11%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.19"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.51%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.47"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -4.94%  "

$ws.Range("E35").Value = "  -0.50%  "

$ws.Range("E36").Value = "  -0.93%  "

$ws.Range("D37").Value = "3.425.59"
$ws.Range("E37").Value = "  -1.64%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.09"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.39%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "167.08"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.27%  "

$ws.Range("E40").Value = "  -1.82%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0785"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.03%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "26.96"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.38%  "

$ws.Range("E43").Value = "  -1.31%  "

$ws.Range("E44").Value = "  +0.12%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.48"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.06%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "41.91"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.07%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.71"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.69%  "

$ws.Range("D48").Value = "2.538.19"
$ws.Range("E48").Value = "  -2.70%  "

$ws.Range("E49").Value = "  -3.50%  "

$ws.Range("E50").Value = "  +0.20%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "22.98"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.58%  "
